# The "endereco" (address) and "cidade" (city) columns are dropped to
# shrink the file for Power BI / notebook consumption. This shifts the
# remaining columns (uf, data_abertura, tipo_agencia) left, from
# E/F/G down to C/D/E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("C:D").Delete()
